$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.072.21'
$ws.Range('E2').Value = '  -0.50%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.007.52'
$ws.Range('E3').Value = '  -1.50%  '
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '224.93'
$ws.Range('E5').Value = '  -1.47%  '
$ws.Range('E6').Value = '  -0.61%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '54.85'
$ws.Range('E8').Value = '  -2.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.372'
$ws.Range('E9').Value = '  -2.69%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0775'
$ws.Range('E10').Value = '  -3.68%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.102'
$ws.Range('E11').Value = '  -4.61%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.303.22'
$ws.Range('E12').Value = '  -1.71%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '13.94'
$ws.Range('E13').Value = '  -3.56%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '19.61'
$ws.Range('E14').Value = '  -3.86%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.15'
$ws.Range('E15').Value = '  -1.86%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.731'
$ws.Range('E16').Value = '  -2.36%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.993.95'
$ws.Range('E17').Value = '  -2.20%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.019.37'
$ws.Range('E18').Value = '  -0.31%  '
$ws.Range('E19').Value = '  +3.86%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.18'
$ws.Range('E20').Value = '  -1.76%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0808'
$ws.Range('E21').Value = '  -3.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '223.11'
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.42'
$ws.Range('E24').Value = '  +2.66%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.14'
$ws.Range('E25').Value = '  -4.89%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.07'
$ws.Range('E26').Value = '  -2.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.86'
$ws.Range('E27').Value = '  -6.25%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.52'
$ws.Range('E28').Value = '  -1.69%  '
$ws.Range('E29').Value = '  -3.48%  '
$ws.Range('E30').Value = '  -6.96%  '
$ws.Range('E31').Value = '  -1.19%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.39'
$ws.Range('E32').Value = '  -1.95%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0597'
$ws.Range('E33').Value = '  -1.70%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.46'
$ws.Range('E34').Value = '  -1.78%  '
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.86'
$ws.Range('E35').Value = '  +2.42%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.30'
$ws.Range('E36').Value = '  -3.09%  '
$ws.Range('E37').Value = '  +0.29%  '
$ws.Range('E38').Value = '  -1.94%  '
$ws.Range('E39').Value = '  -0.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.453.34'
$ws.Range('E40').Value = '  -1.90%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '94.40'
$ws.Range('E41').Value = '  -0.71%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0210'
$ws.Range('E42').Value = '  -3.80%  '
$ws.Range('B43').Value = 'Cronos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0906'
$ws.Range('E43').Value = '  -3.14%  '
$ws.Range('B44').Value = 'HuobiToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.75'
$ws.Range('E44').Value = '  -4.87%  '
$ws.Range('B45').Value = 'FTXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.11'
$ws.Range('E45').Value = '  +15.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '15.81'
$ws.Range('E46').Value = '  -5.56%  '
$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.11'
$ws.Range('E47').Value = '  -2.15%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.992'
$ws.Range('E48').Value = '  -1.42%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.05'
$ws.Range('E49').Value = '  -0.48%  '
$ws.Range('E50').Value = '  -0.82%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.192.29'
$ws.Range('E51').Value = '  -1.75%  '
